$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Correct the floating-point rounding of the timestamp written for the
#    previous batch of rows (828-841): 44232.91759111339 -> 44232.91759111111
# ---------------------------------------------------------------------------
for ($r = 828; $r -le 841; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.91759111111
}

# ---------------------------------------------------------------------------
# 2) Append a new batch of 14 monitored-service rows (842-855), mirroring the
#    same 14-row pattern used throughout the sheet (one row per service).
# ---------------------------------------------------------------------------
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")

# Text shown in the cell (column B) - reuses the same shared strings already
# present in the workbook.
$displayUrls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

# Underlying hyperlink target (Address) - split from SubAddress for the
# MapStore row, same as every other occurrence of this row type.
$addresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$subAddresses = @($null,$null,$null,$null,$null,$null,$null,$null,"/",$null,$null,$null,$null,$null)

$timestamp = 44232.93881568003
$startRow = 842

for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $names[$i]

    $ws.Cells.Item($r, 2).Value = $displayUrls[$i]
    $ws.Cells.Item($r, 2).Style = "Hyperlink"
    if ($subAddresses[$i]) {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $addresses[$i], $subAddresses[$i])
    } else {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $addresses[$i])
    }

    $ws.Cells.Item($r, 3).Value = "Disponible"

    $ws.Cells.Item($r, 4).Value = $timestamp
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
